$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 125 (ALC)
$ws.Range("H125").Value = 18396.95
$ws.Range("I125").Value = 66740
$ws.Range("J125").Value = 2282.6
$ws.Range("K125").Value = 600660
$ws.Range("L125").Value = 20543.4
$ws.Range("M125").Value = -598200
$ws.Range("N125").Value = -25463.4

# Row 132 (ALC)
$ws.Range("H132").Value = 4196.3125
$ws.Range("I132").Value = 3367.1428
$ws.Range("J132").Value = 10000.5
$ws.Range("K132").Value = 10101.4284
$ws.Range("L132").Value = 30001.5
$ws.Range("M132").Value = -7571.428400000001
$ws.Range("N132").Value = -35061.5

# Row 138 (ALC)
$ws.Range("H138").Value = 3022.8948
$ws.Range("I138").Value = 5555
$ws.Range("J138").Value = 2954.4595
$ws.Range("K138").Value = 16665
$ws.Range("L138").Value = 8863.378499999999
$ws.Range("M138").Value = -11525
$ws.Range("N138").Value = -19143.3785

$ws = $wb.Worksheets.Item("ARM")
# Row 74 (ARM)
$ws.Range("H74").Value = 1290.7222
$ws.Range("I74").Value = 1368.8667
$ws.Range("J74").Value = 900
$ws.Range("K74").Value = 1368.8667
$ws.Range("L74").Value = 900
$ws.Range("M74").Value = -494.8667
$ws.Range("N74").Value = -2648

# Row 77 (ARM)
$ws.Range("H77").Value = 1290.7222
$ws.Range("I77").Value = 1368.8667
$ws.Range("J77").Value = 900
$ws.Range("K77").Value = 6844.333500000001
$ws.Range("L77").Value = 4500
$ws.Range("M77").Value = -2476.333500000001
$ws.Range("N77").Value = -13236

# Row 122 (ARM)
$ws.Range("H122").Value = 1595.7878
$ws.Range("J122").Value = 2934.3333
$ws.Range("L122").Value = 8802.999899999999
$ws.Range("N122").Value = -13702.9999

$ws = $wb.Worksheets.Item("BSM")
# Row 105 (BSM)
$ws.Range("H105").Value = 6192.5293
$ws.Range("I105").Value = 6775.143
$ws.Range("J105").Value = 3473.6667
$ws.Range("K105").Value = 6775.143
$ws.Range("L105").Value = 3473.6667
$ws.Range("M105").Value = -5028.143
$ws.Range("N105").Value = -6967.6667

# Row 134 (BSM)
$ws.Range("H134").Value = 1566.0217
$ws.Range("I134").Value = 1350.95
$ws.Range("J134").Value = 2999.8333
$ws.Range("K134").Value = 4052.85
$ws.Range("L134").Value = 8999.499899999999
$ws.Range("M134").Value = -1517.85
$ws.Range("N134").Value = -14069.4999

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 1837.4166
$ws.Range("I31").Value = 1428.7646
$ws.Range("J31").Value = 2829.8572
$ws.Range("K31").Value = 1428.7646
$ws.Range("L31").Value = 2829.8572
$ws.Range("M31").Value = -1133.7646
$ws.Range("N31").Value = -3419.8572

# Row 34 (CRP)
$ws.Range("H34").Value = 1837.4166
$ws.Range("I34").Value = 1428.7646
$ws.Range("J34").Value = 2829.8572
$ws.Range("K34").Value = 1428.7646
$ws.Range("L34").Value = 2829.8572
$ws.Range("M34").Value = -1226.7646
$ws.Range("N34").Value = -3233.8572

# Row 62 (CRP)
$ws.Range("H62").Value = 7002.4443
$ws.Range("I62").Value = 6499.3335
$ws.Range("J62").Value = 8008.6665
$ws.Range("K62").Value = 6499.3335
$ws.Range("L62").Value = 8008.6665
$ws.Range("M62").Value = -5875.3335
$ws.Range("N62").Value = -9256.666499999999

# Row 65 (CRP)
$ws.Range("H65").Value = 7002.4443
$ws.Range("I65").Value = 6499.3335
$ws.Range("J65").Value = 8008.6665
$ws.Range("K65").Value = 32496.6675
$ws.Range("L65").Value = 40043.3325
$ws.Range("M65").Value = -29376.6675
$ws.Range("N65").Value = -46283.3325

# Row 99 (CRP)
$ws.Range("H99").Value = 19350.055
$ws.Range("J99").Value = 8000.1816
$ws.Range("L99").Value = 8000.1816
$ws.Range("N99").Value = -10996.1816

# Row 126 (CRP)
$ws.Range("H126").Value = 19350.055
$ws.Range("J126").Value = 8000.1816
$ws.Range("L126").Value = 24000.5448
$ws.Range("N126").Value = -28940.5448

# Row 132 (CRP)
$ws.Range("H132").Value = 2842.7334
$ws.Range("I132").Value = 2865.9546
$ws.Range("J132").Value = 2778.875
$ws.Range("K132").Value = 8597.863799999999
$ws.Range("L132").Value = 8336.625
$ws.Range("M132").Value = -6067.863799999999
$ws.Range("N132").Value = -13396.625

# Row 134 (CRP)
$ws.Range("H134").Value = 2504.8484
$ws.Range("I134").Value = 2444.516
$ws.Range("J134").Value = 3440
$ws.Range("K134").Value = 7333.548000000001
$ws.Range("L134").Value = 10320
$ws.Range("M134").Value = -4798.548000000001
$ws.Range("N134").Value = -15390

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (CUL)
$ws.Range("H5").Value = 825.5714
$ws.Range("I5").Value = 754.8
$ws.Range("J5").Value = 1002.5
$ws.Range("K5").Value = 2264.4
$ws.Range("L5").Value = 3007.5
$ws.Range("M5").Value = -2152.4
$ws.Range("N5").Value = -3231.5

# Row 37 (CUL)
$ws.Range("H37").Value = 53333.332
$ws.Range("J37").Value = 53333.332
$ws.Range("L37").Value = 159999.996
$ws.Range("N37").Value = -160223.996

# Row 68 (CUL)
$ws.Range("H68").Value = 1666.3334
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 1999
$ws.Range("K68").Value = 4500
$ws.Range("L68").Value = 5997
$ws.Range("M68").Value = -3689
$ws.Range("N68").Value = -7619

# Row 71 (CUL)
$ws.Range("H71").Value = 1666.3334
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 1999
$ws.Range("K71").Value = 13500
$ws.Range("L71").Value = 17991
$ws.Range("M71").Value = -9444
$ws.Range("N71").Value = -26103

# Row 92 (CUL)
$ws.Range("H92").Value = 868
$ws.Range("I92").Value = 882.125
$ws.Range("J92").Value = 845.4
$ws.Range("K92").Value = 2646.375
$ws.Range("L92").Value = 2536.2
$ws.Range("M92").Value = -1398.375
$ws.Range("N92").Value = -5032.2

# Row 113 (CUL)
$ws.Range("H113").Value = 1191.1818
$ws.Range("I113").Value = 1474.75
$ws.Range("J113").Value = 1029.1428
$ws.Range("K113").Value = 4424.25
$ws.Range("L113").Value = 3087.4284
$ws.Range("M113").Value = -2254.25
$ws.Range("N113").Value = -7427.428400000001

# Row 135 (CUL)
$ws.Range("H135").Value = 825.5714
$ws.Range("I135").Value = 754.8
$ws.Range("J135").Value = 1002.5
$ws.Range("K135").Value = 6793.2
$ws.Range("L135").Value = 9022.5
$ws.Range("M135").Value = -4258.2
$ws.Range("N135").Value = -14092.5

# Row 137 (CUL)
$ws.Range("H137").Value = 5560778
$ws.Range("I137").Value = 14287608
$ws.Range("J137").Value = 7341.1816
$ws.Range("K137").Value = 42862824
$ws.Range("L137").Value = 22023.5448
$ws.Range("M137").Value = -42857724
$ws.Range("N137").Value = -32223.5448

$ws = $wb.Worksheets.Item("GSM")
# Row 122 (GSM)
$ws.Range("H122").Value = 2600.4
$ws.Range("I122").Value = 976
$ws.Range("J122").Value = 3683.3333
$ws.Range("K122").Value = 2928
$ws.Range("L122").Value = 11049.9999
$ws.Range("M122").Value = -478
$ws.Range("N122").Value = -15949.9999

# Row 126 (GSM)
$ws.Range("H126").Value = 6822.9443
$ws.Range("I126").Value = 7186.8
$ws.Range("K126").Value = 21560.4
$ws.Range("M126").Value = -19090.4

# Row 132 (GSM)
$ws.Range("H132").Value = 3175
$ws.Range("I132").Value = 3356.8333
$ws.Range("J132").Value = 1538.5
$ws.Range("K132").Value = 10070.4999
$ws.Range("L132").Value = 4615.5
$ws.Range("M132").Value = -7540.499899999999
$ws.Range("N132").Value = -9675.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Range("H7").Value = 22728.625
$ws.Range("I7").Value = 29423.455
$ws.Range("K7").Value = 29423.455
$ws.Range("M7").Value = -29311.455

# Row 16 (LTW)
$ws.Range("H16").Value = 1152.5238
$ws.Range("I16").Value = 1206.375
$ws.Range("J16").Value = 980.2
$ws.Range("K16").Value = 1206.375
$ws.Range("L16").Value = 980.2
$ws.Range("M16").Value = -1036.375
$ws.Range("N16").Value = -1320.2

# Row 126 (LTW)
$ws.Range("H126").Value = 22728.625
$ws.Range("I126").Value = 29423.455
$ws.Range("K126").Value = 88270.36500000001
$ws.Range("M126").Value = -85800.36500000001

# Row 132 (LTW)
$ws.Range("H132").Value = 2664.8704
$ws.Range("I132").Value = 2134.675
$ws.Range("J132").Value = 4179.7144
$ws.Range("K132").Value = 6404.025000000001
$ws.Range("L132").Value = 12539.1432
$ws.Range("M132").Value = -3874.025000000001
$ws.Range("N132").Value = -17599.1432

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (WVR)
$ws.Range("H132").Value = 8522.380999999999
$ws.Range("I132").Value = 11287.333
$ws.Range("J132").Value = 1610
$ws.Range("K132").Value = 33861.999
$ws.Range("L132").Value = 4830
$ws.Range("M132").Value = -31331.999
$ws.Range("N132").Value = -9890

# Row 136 (WVR)
$ws.Range("H136").Value = 933.0625
$ws.Range("I136").Value = 661.93335
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 1985.80005
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = 564.1999499999999
$ws.Range("N136").Value = -20100
